# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# D-column price cells are forced to text with a leading apostrophe so Excel
# doesn't reinterpret dotted price strings (e.g. "3.318.13") as numbers/dates;
# E-column percentage cells already contain spaces/"%" so they stay text as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.860.18"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "'3.318.13"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'557.82"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'186.01"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'3.311.72"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("E10").Value = "  -6.54%  "
$ws.Range("D11").Value = "'0.578"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "'45.77"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "'3.847.30"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "'570.82"
$ws.Range("E16").Value = "  -9.31%  "
$ws.Range("D17").Value = "'65.782.86"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.117"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.318.96"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "'10.83"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("D22").Value = "'0.889"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").Value = "'17.95"
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").Value = "'5.00"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "'97.65"
$ws.Range("E25").Value = "  -8.41%  "
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Value = "'30.44"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'6.73"
$ws.Range("E31").Value = "  +7.47%  "
$ws.Range("D32").Value = "'3.69"
$ws.Range("E32").Value = "  -8.48%  "
$ws.Range("D33").Value = "'559.32"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("D34").Value = "'10.83"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("D36").Value = "'3.746.83"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'55.47"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").Value = "'33.74"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("E40").Value = "  -4.04%  "
$ws.Range("E41").Value = "  -5.75%  "
$ws.Range("D42").Value = "'2.58"
$ws.Range("E42").Value = "  -5.76%  "
$ws.Range("E43").Value = "  -8.32%  "
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "'2.98"
$ws.Range("E47").Value = "  -12.48%  "
$ws.Range("D48").Value = "'0.127"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").Value = "'125.09"
$ws.Range("E51").Value = "  +2.47%  "
